$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.721.78'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '1.864.00'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.012'
$ws.Range("E4").Value = '  +0.92%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.27'
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.011'
$ws.Range("E6").Value = '  +0.90%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4673'
$ws.Range("E7").Value = '  -1.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3896'
$ws.Range("E8").Value = '  -1.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.64'
$ws.Range("E9").Value = '  -2.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07992'
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.000'
$ws.Range("E11").Value = '  -2.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.55'
$ws.Range("E12").Value = '  -2.62%  '
$ws.Range("D13").Value = '1.867.94'
$ws.Range("E13").Value = '  -1.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.983'
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.157'
$ws.Range("E15").Value = '  +0.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.012'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.07'
$ws.Range("E17").Value = '  +1.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06688'
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001042'
$ws.Range("E19").Value = '  -0.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.89'
$ws.Range("E20").Value = '  -2.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.011'
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").Value = '27.709.47'
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.462'
$ws.Range("E23").Value = '  -1.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.89'
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.321'
$ws.Range("E25").Value = '  +0.67%  '
$ws.Range("D26").Value = '2.089.46'
$ws.Range("E26").Value = '  -1.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.05'
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.65'
$ws.Range("E28").Value = '  -2.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.103'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.394'
$ws.Range("E30").Value = '  -3.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '120.89'
$ws.Range("E31").Value = '  -1.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9710'
$ws.Range("E32").Value = '  -0.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09431'
$ws.Range("E33").Value = '  -1.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.641'
$ws.Range("E34").Value = '  +1.52%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.305'
$ws.Range("E35").Value = '  -0.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.335'
$ws.Range("E36").Value = '  -8.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06034'
$ws.Range("E37").Value = '  -1.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02212'
$ws.Range("E38").Value = '  -1.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.202'
$ws.Range("E39").Value = '  -2.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.163'
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.011'
$ws.Range("E41").Value = '  +1.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5911'
$ws.Range("E42").Value = '  -1.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1880'
$ws.Range("E43").Value = '  -1.18%  '
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.258'
$ws.Range("E45").Value = '  +0.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5622'
$ws.Range("E46").Value = '  -1.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.03'
$ws.Range("E47").Value = '  -1.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.917'
$ws.Range("E48").Value = '  -1.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.296'
$ws.Range("E49").Value = '  -2.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06765'
$ws.Range("E50").Value = '  -1.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '112.53'
$ws.Range("E51").Value = '  -2.13%  '
